# ---------------------------------------------------------------------------
# CARLA dataset spec workbook update
#   * Add a new "Dataset_info" sheet (placed between "Specifications" and
#     "Information") summarising number of images per town / phase.
#   * Hide the old "Information" sheet (its data stays intact).
#   * Tweak a handful of cell alignments on the "Information" sheet (the
#     per-scenario "Incident Types" column becomes vertically centred too).
#   * Refresh the selections so the workbook reopens where the author left
#     it.
#
# NOTE: worksheet variables in this host resolve by *position*, not stable
# identity -- once a sheet is inserted/removed, any previously-captured
# reference can silently start pointing at whatever now sits at that index.
# So every sheet handle used after a structural change (Add/Move/Delete) is
# re-fetched by name right before it's used.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "Dataset_info" worksheet right before "Information".
# ---------------------------------------------------------------------------
$wsInfoForInsert = $wb.Worksheets.Item("Information")
$ds = $wb.Worksheets.Add($wsInfoForInsert)
$ds.Name = "Dataset_info"

# Header row
$ds.Range("A1").Value2 = "Town Map"
$ds.Range("B1").Value2 = "Phases"
$ds.Range("C1").Value2 = "No of images"

# ---------------------------------------------------------------------------
# Per-town / per-phase image counts. $null means "value left blank" (the
# shoot for that phase has not happened yet at the time of this snapshot).
# ---------------------------------------------------------------------------
$towns = @(
    @{ Name = "Town_01"; StartRow = 2;  Images = @(100, 200, 100, 100, 200, 100) },
    @{ Name = "Town_02"; StartRow = 8;  Images = @(300, 100, 200) },
    @{ Name = "Town_03"; StartRow = 11; Images = @(100, 200, 100, 100, 200, 100, 200) },
    @{ Name = "Town_04"; StartRow = 18; Images = @(670, $null, $null, $null, $null) },
    @{ Name = "Town_05"; StartRow = 23; Images = @($null, $null, $null, $null) },
    @{ Name = "Town_06"; StartRow = 27; Images = @($null, $null, $null) },
    @{ Name = "Town_07"; StartRow = 30; Images = @($null, $null, $null) }
)

foreach ($town in $towns) {
    $startRow = $town.StartRow
    $count = $town.Images.Count
    $endRow = $startRow + $count - 1

    for ($i = 0; $i -lt $count; $i++) {
        $row = $startRow + $i
        $ds.Range("B$row").Value2 = "Phase " + ($i + 1)
        if ($null -ne $town.Images[$i]) {
            $ds.Range("C$row").Value2 = $town.Images[$i]
        }
    }

    $ds.Range("A${startRow}:A$endRow").Value2 = $town.Name
    if ($endRow -gt $startRow) {
        $ds.Range("A${startRow}:A$endRow").Merge()
    }
}

# Totals row
$ds.Range("A33").Value2 = "Total images"
$ds.Range("A33:B33").Merge()
$ds.Range("C33").Formula = "=SUM(C2:C32)"

# ---------------------------------------------------------------------------
# Formatting: every populated data cell gets a thin box border; column B/C
# values are centred, column A town labels are centred (+ vertically
# centred for the tallest merged block).
# ---------------------------------------------------------------------------
$ds.Range("A1:C33").Borders.LineStyle = 1
$ds.Range("A1:C33").HorizontalAlignment = -4108

$ds.Range("A2:A32").HorizontalAlignment = -4108
$ds.Range("A11:A17").VerticalAlignment = -4108

$ds.Range("A33:B33").HorizontalAlignment = -4108

$ds.Columns.Item(1).AutoFit() | Out-Null
$ds.Columns.Item(3).AutoFit() | Out-Null

$ds.Range("E25").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Hide the old "Information" sheet (content untouched).
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Information")
$wsInfo.Visible = $false

# ---------------------------------------------------------------------------
# 3. On "Information", the first "Incident Types" cell of each scenario
#    block (merged over 3 rows) becomes vertically centred as well.
# ---------------------------------------------------------------------------
$wsInfo2 = $wb.Worksheets.Item("Information")
$incidentBlocks = @("B2:B4", "B14:B16", "B26:B28", "B38:B40", "B53:B55")
foreach ($rng in $incidentBlocks) {
    $wsInfo2.Range($rng).VerticalAlignment = -4108
}
$wsInfo2.Range("B11:B13").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Restore the "Specifications" selection recorded by the author.
# ---------------------------------------------------------------------------
$wsSpec = $wb.Worksheets.Item("Specifications")
$wsSpec.Range("B24").Select() | Out-Null

$ds2 = $wb.Worksheets.Item("Dataset_info")
$ds2.Activate() | Out-Null
